$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 32 (ALC)
$ws.Range("H32").Value = 610
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 610
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 610
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -1262

# Row 34 (ALC)
$ws.Range("H34").Value = 12022.4
$ws.Range("I34").Value = 2528
$ws.Range("J34").Value = 50000
$ws.Range("K34").Value = 2528
$ws.Range("L34").Value = 50000
$ws.Range("M34").Value = -2325

# Row 36 (ALC)
$ws.Range("H36").Value = 12022.4
$ws.Range("I36").Value = 2528
$ws.Range("J36").Value = 50000
$ws.Range("K36").Value = 2528
$ws.Range("L36").Value = 50000
$ws.Range("M36").Value = -1813

# Row 39 (ALC)
$ws.Range("H39").Value = 215
$ws.Range("I39").Value = 95
$ws.Range("J39").Value = 269.54544
$ws.Range("K39").Value = 285
$ws.Range("L39").Value = 808.63632
$ws.Range("M39").Value = 11
$ws.Range("N39").Value = -1400.63632

# Row 43 (ALC)
$ws.Range("H43").Value = 673.5714
$ws.Range("I43").Value = 696
$ws.Range("J43").Value = 661.1111
$ws.Range("K43").Value = 696
$ws.Range("L43").Value = 661.1111
$ws.Range("M43").Value = -627
$ws.Range("N43").Value = -799.1111

# Row 76 (ALC)
$ws.Range("H76").Value = 3271231.8
$ws.Range("I76").Value = 4118336.5
$ws.Range("J76").Value = 3827.1428
$ws.Range("K76").Value = 4118336.5
$ws.Range("L76").Value = 3827.1428
$ws.Range("M76").Value = -4118021.5

# Row 79 (ALC)
$ws.Range("H79").Value = 3271231.8
$ws.Range("I79").Value = 4118336.5
$ws.Range("J79").Value = 3827.1428
$ws.Range("K79").Value = 4118336.5
$ws.Range("L79").Value = 3827.1428
$ws.Range("M79").Value = -4117244.5

# Row 134 (ALC)
$ws.Range("H134").Value = 65431.43
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 65431.43
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 65431.43
$ws.Range("N134").Value = -75571.42999999999

# Row 137 (ALC)
$ws.Range("H137").Value = 20834562
$ws.Range("I137").Value = 27778778
$ws.Range("J137").Value = 1912.5
$ws.Range("K137").Value = 83336334
$ws.Range("L137").Value = 5737.5
$ws.Range("M137").Value = -83333784
$ws.Range("N137").Value = -10837.5

$ws = $wb.Worksheets.Item("ARM")
# Row 4 (ARM)
$ws.Range("H4").Value = 175
$ws.Range("I4").Value = 175
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 175
$ws.Range("L4").Value = 0
$ws.Range("M4").Value = -59
$ws.Range("N4").ClearContents()

# Row 23 (ARM)
$ws.Range("H23").Value = 19000
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 19000
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 19000
$ws.Range("M23").ClearContents()
$ws.Range("N23").Value = -19518

# Row 32 (ARM)
$ws.Range("H32").Value = 21748.926
$ws.Range("I32").Value = 4261.49
$ws.Range("J32").Value = 77490.125
$ws.Range("K32").Value = 4261.49
$ws.Range("L32").Value = 77490.125
$ws.Range("M32").Value = -3974.49
$ws.Range("N32").Value = -78064.125

# Row 37 (ARM)
$ws.Range("H37").Value = 4300
$ws.Range("I37").Value = 0
$ws.Range("J37").Value = 4300
$ws.Range("K37").Value = 0
$ws.Range("L37").Value = 4300
$ws.Range("N37").Value = -4846

# Row 63 (ARM)
$ws.Range("H63").Value = 4120.5
$ws.Range("I63").Value = 4299.737
$ws.Range("J63").Value = 3742.111
$ws.Range("K63").Value = 4299.737
$ws.Range("L63").Value = 3742.111
$ws.Range("M63").Value = -3613.737
$ws.Range("N63").Value = -5114.111

# Row 66 (ARM)
$ws.Range("H66").Value = 4120.5
$ws.Range("I66").Value = 4299.737
$ws.Range("J66").Value = 3742.111
$ws.Range("K66").Value = 21498.685
$ws.Range("L66").Value = 18710.555
$ws.Range("M66").Value = -18066.685
$ws.Range("N66").Value = -25574.555

# Row 74 (ARM)
$ws.Range("H74").Value = 7531.35
$ws.Range("I74").Value = 1053.7693
$ws.Range("J74").Value = 19561.143
$ws.Range("K74").Value = 1053.7693
$ws.Range("L74").Value = 19561.143
$ws.Range("M74").Value = -179.7692999999999
$ws.Range("N74").Value = -21309.143

# Row 77 (ARM)
$ws.Range("H77").Value = 7531.35
$ws.Range("I77").Value = 1053.7693
$ws.Range("J77").Value = 19561.143
$ws.Range("K77").Value = 5268.8465
$ws.Range("L77").Value = 97805.715
$ws.Range("M77").Value = -900.8464999999997
$ws.Range("N77").Value = -106541.715

# Row 122 (ARM)
$ws.Range("H122").Value = 6849.615
$ws.Range("I122").Value = 7909.2
$ws.Range("J122").Value = 6187.375
$ws.Range("K122").Value = 23727.6
$ws.Range("L122").Value = 18562.125
$ws.Range("M122").Value = -21277.6
$ws.Range("N122").Value = -23462.125

# Row 132 (ARM)
$ws.Range("H132").Value = 3251.1072
$ws.Range("I132").Value = 2877.8696
$ws.Range("J132").Value = 4968
$ws.Range("K132").Value = 8633.6088
$ws.Range("L132").Value = 14904
$ws.Range("M132").Value = -6103.6088
$ws.Range("N132").Value = -19964

$ws = $wb.Worksheets.Item("BSM")
# Row 107 (BSM)
$ws.Range("H107").Value = 750.9091
$ws.Range("I107").Value = 494.2857
$ws.Range("J107").Value = 1200
$ws.Range("K107").Value = 494.2857
$ws.Range("L107").Value = 1200
$ws.Range("M107").Value = 1425.7143
$ws.Range("N107").Value = -5040

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (CRP)
$ws.Range("H31").Value = 5454.135
$ws.Range("I31").Value = 1842.4348
$ws.Range("J31").Value = 8318.585999999999
$ws.Range("K31").Value = 1842.4348
$ws.Range("L31").Value = 8318.585999999999
$ws.Range("M31").Value = -1547.4348
$ws.Range("N31").Value = -8908.585999999999

# Row 34 (CRP)
$ws.Range("H34").Value = 5454.135
$ws.Range("I34").Value = 1842.4348
$ws.Range("J34").Value = 8318.585999999999
$ws.Range("K34").Value = 1842.4348
$ws.Range("L34").Value = 8318.585999999999
$ws.Range("M34").Value = -1640.4348
$ws.Range("N34").Value = -8722.585999999999

# Row 107 (CRP)
$ws.Range("H107").Value = 2278.25
$ws.Range("I107").Value = 2266.6667
$ws.Range("J107").Value = 2313
$ws.Range("K107").Value = 2266.6667
$ws.Range("L107").Value = 2313
$ws.Range("M107").Value = -346.6667000000002
$ws.Range("N107").Value = -6153

# Row 132 (CRP)
$ws.Range("H132").Value = 2707.3076
$ws.Range("I132").Value = 1294.375
$ws.Range("J132").Value = 4968
$ws.Range("K132").Value = 3883.125
$ws.Range("L132").Value = 14904
$ws.Range("M132").Value = -1353.125

$ws = $wb.Worksheets.Item("CUL")
# Row 22 (CUL)
$ws.Range("H22").Value = 1280
$ws.Range("I22").Value = 500
$ws.Range("J22").Value = 1475
$ws.Range("K22").Value = 1500
$ws.Range("L22").Value = 4425
$ws.Range("M22").Value = -1331
$ws.Range("N22").Value = -4763

# Row 27 (CUL)
$ws.Range("H27").Value = 1280
$ws.Range("I27").Value = 500
$ws.Range("J27").Value = 1475
$ws.Range("K27").Value = 1500
$ws.Range("L27").Value = 4425
$ws.Range("M27").Value = -1398
$ws.Range("N27").Value = -4629

# Row 58 (CUL)
$ws.Range("H58").Value = 10000
$ws.Range("I58").Value = 0
$ws.Range("J58").Value = 10000
$ws.Range("K58").Value = 0
$ws.Range("L58").Value = 30000
$ws.Range("M58").ClearContents()
$ws.Range("N58").Value = -30256

# Row 86 (CUL)
$ws.Range("H86").Value = 3000
$ws.Range("I86").Value = 0
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 0
$ws.Range("L86").Value = 9000
$ws.Range("M86").ClearContents()
$ws.Range("N86").Value = -11372

# Row 89 (CUL)
$ws.Range("H89").Value = 3000
$ws.Range("I89").Value = 0
$ws.Range("J89").Value = 3000
$ws.Range("K89").Value = 0
$ws.Range("L89").Value = 27000
$ws.Range("M89").ClearContents()
$ws.Range("N89").Value = -38856

# Row 107 (CUL)
$ws.Range("H107").Value = 326.33334
$ws.Range("I107").Value = 334.0476
$ws.Range("J107").Value = 319.58334
$ws.Range("K107").Value = 1002.1428
$ws.Range("L107").Value = 958.7500200000001
$ws.Range("M107").Value = 917.8572

# Row 113 (CUL)
$ws.Range("H113").Value = 27778308
$ws.Range("I113").Value = 400
$ws.Range("J113").Value = 29412302
$ws.Range("K113").Value = 1200
$ws.Range("L113").Value = 88236906
$ws.Range("M113").Value = 970

# Row 131 (CUL)
$ws.Range("H131").Value = 6668149
$ws.Range("I131").Value = 567.1429000000001
$ws.Range("J131").Value = 7753569.5
$ws.Range("K131").Value = 1701.4287
$ws.Range("L131").Value = 23260708.5
$ws.Range("M131").Value = 3338.5713
$ws.Range("N131").Value = -23270788.5

# Row 134 (CUL)
$ws.Range("H134").Value = 6145.1763
$ws.Range("I134").Value = 2736.9
$ws.Range("J134").Value = 11014.143
$ws.Range("K134").Value = 8210.700000000001
$ws.Range("L134").Value = 33042.429
$ws.Range("M134").Value = -3140.700000000001
$ws.Range("N134").Value = -43182.429

# Row 140 (CUL)
$ws.Range("H140").Value = 7419.722
$ws.Range("I140").Value = 10277.728
$ws.Range("J140").Value = 2928.5715
$ws.Range("K140").Value = 30833.184
$ws.Range("L140").Value = 8785.7145
$ws.Range("M140").Value = -25653.184

$ws = $wb.Worksheets.Item("GSM")
# Row 122 (GSM)
$ws.Range("H122").Value = 2081.1428
$ws.Range("I122").Value = 1876
$ws.Range("J122").Value = 2833.3333
$ws.Range("K122").Value = 5628
$ws.Range("L122").Value = 8499.999899999999
$ws.Range("M122").Value = -3178

# Row 123 (GSM)
$ws.Range("H123").Value = 10484.52
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 10484.52
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 10484.52
$ws.Range("N123").Value = -15384.52

# Row 126 (GSM)
$ws.Range("H126").Value = 2382.2334
$ws.Range("I126").Value = 1968.4445
$ws.Range("J126").Value = 2559.5715
$ws.Range("K126").Value = 5905.333500000001
$ws.Range("L126").Value = 7678.7145
$ws.Range("M126").Value = -3435.333500000001
$ws.Range("N126").Value = -12618.7145

# Row 132 (GSM)
$ws.Range("H132").Value = 4421.1665
$ws.Range("I132").Value = 4877.4287
$ws.Range("J132").Value = 3782.4
$ws.Range("K132").Value = 14632.2861
$ws.Range("L132").Value = 11347.2
$ws.Range("M132").Value = -12102.2861
$ws.Range("N132").Value = -16407.2

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (LTW)
$ws.Range("H7").Value = 3220.25
$ws.Range("I7").Value = 2400
$ws.Range("J7").Value = 3425.3125
$ws.Range("K7").Value = 2400
$ws.Range("L7").Value = 3425.3125
$ws.Range("M7").Value = -2288
$ws.Range("N7").Value = -3649.3125

# Row 22 (LTW)
$ws.Range("H22").Value = 1013.5714
$ws.Range("I22").Value = 786.6667
$ws.Range("J22").Value = 1183.75
$ws.Range("K22").Value = 786.6667
$ws.Range("L22").Value = 1183.75
$ws.Range("M22").Value = -491.6667
$ws.Range("N22").Value = -1773.75

# Row 27 (LTW)
$ws.Range("H27").Value = 1013.5714
$ws.Range("I27").Value = 786.6667
$ws.Range("J27").Value = 1183.75
$ws.Range("K27").Value = 786.6667
$ws.Range("L27").Value = 1183.75
$ws.Range("M27").Value = -679.6667
$ws.Range("N27").Value = -1397.75

# Row 126 (LTW)
$ws.Range("H126").Value = 3220.25
$ws.Range("I126").Value = 2400
$ws.Range("J126").Value = 3425.3125
$ws.Range("K126").Value = 7200
$ws.Range("L126").Value = 10275.9375
$ws.Range("M126").Value = -4730
$ws.Range("N126").Value = -15215.9375

$ws = $wb.Worksheets.Item("WVR")
# Row 81 (WVR)
$ws.Range("H81").Value = 3467.1538
$ws.Range("I81").Value = 1316.2222
$ws.Range("J81").Value = 4605.8823
$ws.Range("K81").Value = 2632.4444
$ws.Range("L81").Value = 9211.7646
$ws.Range("M81").Value = -1571.4444
$ws.Range("N81").Value = -11333.7646

# Row 84 (WVR)
$ws.Range("H84").Value = 3467.1538
$ws.Range("I84").Value = 1316.2222
$ws.Range("J84").Value = 4605.8823
$ws.Range("K84").Value = 13162.222
$ws.Range("L84").Value = 46058.823
$ws.Range("M84").Value = -7858.222
$ws.Range("N84").Value = -56666.823

# Row 122 (WVR)
$ws.Range("H122").Value = 1926
$ws.Range("I122").Value = 1952
$ws.Range("J122").Value = 1900
$ws.Range("K122").Value = 5856
$ws.Range("L122").Value = 5700
$ws.Range("M122").Value = -3406
$ws.Range("N122").Value = -10600

# Row 123 (WVR)
$ws.Range("H123").Value = 32800
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 32800
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 32800
$ws.Range("N123").Value = -42600

# Row 132 (WVR)
$ws.Range("H132").Value = 12999.667
$ws.Range("I132").Value = 21333.334
$ws.Range("J132").Value = 4666
$ws.Range("K132").Value = 64000.00199999999
$ws.Range("L132").Value = 13998
$ws.Range("M132").Value = -61470.00199999999
$ws.Range("N132").Value = -19058

